$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, centered, thin border) from H1 onto the
# two new header cells I1:J1, then set their text values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new I0 / IF data columns for rows 2-14.
$data = @{
    2  = @(6, 6)
    3  = @(9, 9)
    4  = @(10, 10)
    5  = @(8, 8)
    6  = @(7, 7)
    7  = @(8, 8)
    8  = @(5, 6)
    9  = @(1, 1)
    10 = @(9, 9)
    11 = @(9, 9)
    12 = @(5, 5)
    13 = @(4, 4)
    14 = @(7, 7)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("I$row").Value = $vals[0]
    $ws.Range("J$row").Value = $vals[1]
}
